# Fruta / hortaliza, semanal
# New weekly observation is inserted at row 191 (pushing the existing rows
# 191-208 down to 192-209); row 191 is then populated with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(191).Insert()

$ws.Range("A191").Value = 10
$ws.Range("B191").Value = "Vega Modelo de Temuco"
$ws.Range("C191").Value = "La Araucanía"
$ws.Range("D191").Value = 44461
$ws.Range("E191").Value = 9
$ws.Range("F191").Value = "Fruta"
$ws.Range("G191").Value = 100108
$ws.Range("H191").Value = "Tropicales y subtropicales"
$ws.Range("I191").Value = 100108002
$ws.Range("J191").Value = "Mango"
$ws.Range("K191").Value = "Sin especificar"
$ws.Range("L191").Value = "Primera"
$ws.Range("M191").Value = 500
$ws.Range("N191").Value = 9000
$ws.Range("O191").Value = 9000
$ws.Range("P191").Value = 9000
$ws.Range("Q191").Value = "$/bandeja 4 kilos"
$ws.Range("R191").Value = "Brasil"
$ws.Range("S191").Value = 2250
$ws.Range("T191").Value = 4
